# Daily attendance processing - 2025-11-24 05:53:04
# Reorders the "Recorded By" column (G) entries so that when the value is
# exactly "System, <email>" it becomes "<email>, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.ToString().StartsWith("System, ")) {
        $rest = $val.ToString().Substring(8)
        $cell.Value = "$rest, System"
    }
}
